# The sheet originally listed years 2001..2019 in rows 2..20 (row 1 is the
# header). The data was refreshed to cover 2010..2020 instead: the first nine
# data rows (2001年..2009年, rows 2-10) are dropped, which shifts the former
# 2010年..2019年 rows (11-20) up to become rows 2-11, and a brand-new 2020年
# row is appended (row 12) with its own 其他(净) / 贷款 / 资金来源 figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the nine oldest years (2001年-2009年). Everything below slides up,
# so the old "2010年" row (11) becomes the new row 2, ... old "2019年" (20)
# becomes the new row 11.
$ws.Rows("2:10").Delete()

# Append the new 2020年 row, cloning the formatting (bold/centered/bordered
# year label, plain number cells, etc.) from the row directly above it
# (now row 11, 2019年) so the new row matches the look of the rest of the
# table.
$ws.Range("A11:U11").Copy()
$ws.Range("A12:U12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A12").Value = "2020年"
$ws.Range("D12").Value = 715.846331776013
$ws.Range("Q12").Value = 86052.9493347746
$ws.Range("R12").Value = 86768.79566655061
